# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    3, 2, 0, 1, 2, 1, 2, 3, 3, 1,
    1, 2, 0, 2, 2, 1, 3, 3, 1, 1,
    1, 3, 1, 1, 1, 1, 4, 1, 2, 0,
    2, 1, 2, 3, 1, 3, 0, 1, 2, 1,
    2, 5, 2, 1, 3, 1, 2, 2, 0, 1,
    0, 2, 2, 2, 1, 0, 2, 1, 3, 3,
    2, 5, 3, 1, 1, 2, 3, 2
)

$startRow = 2
for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
